# Update column F (dSF) values for specific rows based on a data repull /
# recalculation of the mean ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = -5
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -4
$ws.Range("F14").Value = -2
$ws.Range("F16").Value = 3
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = -4
$ws.Range("F26").Value = -5
$ws.Range("F30").Value = 0
$ws.Range("F33").Value = -4
